$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 31: new TA sample entry (Feb 15, 2019 -> serial 43511)
$ws.Range("A31").Value = 43511

$ws.Range("B31").Value = 2210.7950000000001
$ws.Range("D31").Formula = "=100*(B31-C31)/C31"
$ws.Range("F31").Value = "New CRM bottle (opened 02/14)"

# Update the view state to match where the user left off
$ws.Range("B32").Select()
$excel.ActiveWindow.ScrollRow = 28
